$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$elem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet updates ---

# Version: 2.0.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-08-01T06:39:39+00:00 -> 2025-12-19T08:22:07+00:00
$meta.Range("B8").Value = "2025-12-19T08:22:07+00:00"

# Description: expand into multi-line markdown bullet text (also mirrored on Elements!M2)
$newDescription = "`n- **Séjour** : commentaire relatif au séjour.`n- **Événement** : commentaires sur le déroulé de l’évènement.`n- **Évaluation** : commentaire libre sur le contenu ou le résultat de l’évaluation.`n- **Champ évalué** : commentaire spécifique à un item ou sous-item évalué."
$meta.Range("B12").Value = $newDescription
$elem.Range("M2").Value = $newDescription

# New row 22: another Context entry for QuestionnaireResponse, matching formatting of row 21
$meta.Range("A21:B21").Copy()
$meta.Range("A22:B22").PasteSpecial(-4122)
$meta.Range("A22").Value = "Context"
$meta.Range("B22").Value = "element:QuestionnaireResponse"

# --- Elements sheet updates ---

# Short column for Extension.value[x] row: replace generic text with concrete example
$elem.Range("L6").Value = "Exemple de commentaire : Cet évènement a débuté plus tard l’usager était sous la douche à l’heure du début du rendez-vous."
